# Update "Forecast Comparison" sheet with the correct forecast output:
#   - insert a new "Week_Start_Date" column right after "Week" (this shifts
#     ASIN / MyForecast / Amazon * Forecast / Product Title / is_holiday_week
#     from columns B..I to C..J)
#   - rewrite the Week labels to drop the leading zero (W01 -> W1, etc.)
#   - populate the new Week_Start_Date column with the week's start date
#   - correct the forecast figures (MyForecast / Mean / P70 / P80 / P90)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert a new column B ("Week_Start_Date") -----------------------------
$ws.Columns.Item(2).Insert()

# --- Header row --------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Make sure the new date strings ("2025-01-05", ...) are stored as plain
# text rather than being auto-converted to date serial numbers by Excel.
$dateRange = $ws.Range("B2:B17")
$dateRange.NumberFormat = "@"

# --- Row data ---------------------------------------------------------
# Columns after the insert:
#   A=Week, B=Week_Start_Date, C=ASIN, D=MyForecast,
#   E=Amazon Mean Forecast, F=Amazon P70 Forecast, G=Amazon P80 Forecast,
#   H=Amazon P90 Forecast, I=Product Title, J=is_holiday_week
$rows = @(
    @{ Row = 2;  Week = "W1";  Date = "2025-01-05"; MyForecast = 13; Mean = 1; P70 = 1; P80 = 1; P90 = 2 },
    @{ Row = 3;  Week = "W2";  Date = "2025-01-12"; MyForecast = 4;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 4;  Week = "W3";  Date = "2025-01-19"; MyForecast = 6;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 5;  Week = "W4";  Date = "2025-01-26"; MyForecast = 1;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 6;  Week = "W5";  Date = "2025-02-02"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 7;  Week = "W6";  Date = "2025-02-09"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 8;  Week = "W7";  Date = "2025-02-16"; MyForecast = 5;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 9;  Week = "W8";  Date = "2025-02-23"; MyForecast = 1;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 10; Week = "W9";  Date = "2025-03-02"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 4 },
    @{ Row = 11; Week = "W10"; Date = "2025-03-09"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 12; Week = "W11"; Date = "2025-03-16"; MyForecast = 5;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 13; Week = "W12"; Date = "2025-03-23"; MyForecast = 1;  Mean = 1; P70 = 1; P80 = 1; P90 = 2 },
    @{ Row = 14; Week = "W13"; Date = "2025-03-30"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 15; Week = "W14"; Date = "2025-04-06"; MyForecast = 2;  Mean = 2; P70 = 1; P80 = 2; P90 = 3 },
    @{ Row = 16; Week = "W15"; Date = "2025-04-13"; MyForecast = 5;  Mean = 1; P70 = 1; P80 = 1; P90 = 2 },
    @{ Row = 17; Week = "W16"; Date = "2025-04-20"; MyForecast = 1;  Mean = 1; P70 = 1; P80 = 1; P90 = 2 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Week
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 4).Value = $r.MyForecast
    $ws.Cells.Item($row, 5).Value = $r.Mean
    $ws.Cells.Item($row, 6).Value = $r.P70
    $ws.Cells.Item($row, 7).Value = $r.P80
    $ws.Cells.Item($row, 8).Value = $r.P90
}

# Drop the temporary text-number-format override so the cells keep the
# workbook's default (unstyled) look, same as every other text cell.
$dateRange.Style = "Normal"

# The is_holiday_week column (J) keeps its original blank values; clear
# out any leftover content there after the column shift.
$ws.Range("J2:J17").ClearContents()
